$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.715798258781433
$ws.Range("B1").Value = 1.827842831611633
$ws.Range("C1").Value = 1.768250942230225
$ws.Range("D1").Value = 2.114036083221436
$ws.Range("E1").Value = 2.977443695068359
